$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6391
$ws1.Range("F5").Value = 375
$ws1.Range("F8").Value = 23
$ws1.Range("F9").Value = 74
$ws1.Range("F10").Value = 70
$ws1.Range("F13").Value = 364
$ws1.Range("F14").Value = 773
$ws1.Range("F15").Value = 3109
$ws1.Range("F17").Value = 182
$ws1.Range("F18").Value = 1773

# Sheet "全部类型" (all types) - fourth sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6391
$ws4.Range("F5").Value = 375
$ws4.Range("F9").Value = 23
$ws4.Range("F10").Value = 74
$ws4.Range("F11").Value = 70
$ws4.Range("F14").Value = 364
$ws4.Range("F15").Value = 773
$ws4.Range("F16").Value = 3109
$ws4.Range("F18").Value = 182
$ws4.Range("F19").Value = 1773
